$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) Remove the stray "_GoBack" bookmark that currently sits between
#    "...–omogući –PFS" and "- za upravljanje korisicima".
# ------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# ------------------------------------------------------------------
# 2) Extend the "PROBLEMI:Reload" paragraph with the extra remark,
#    re-creating the "_GoBack" bookmark in its new location (between
#    "...prebacuje stra" and "nice,").
# ------------------------------------------------------------------
$r = $d.Content
$found = $r.Find.Execute("PROBLEMI:Reload", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$r.Collapse(0)
$r.InsertAfter(", kod searcha mi ne prebacuje stranice,")

# Re-insert the _GoBack bookmark right after "...prebacuje stra" (i.e.
# immediately before "nice,").
$full = $d.Content
$full.Find.Execute("stra", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$bmPos = $full.End
$d.Bookmarks.Add("_GoBack", $d.Range($bmPos, $bmPos))

# Split "PROBLEMI:Reload" and ", " / "kod searcha..." into separate
# runs, the way the original edit left them, by briefly dropping a
# bookmark at each boundary and removing it again (the removal leaves
# the run split behind while not leaving extra bookmarks around).
$p1 = $d.Content
$p1.Find.Execute("PROBLEMI:Reload", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$split1 = $p1.End
$d.Bookmarks.Add("TmpSplit1", $d.Range($split1, $split1))
$d.Bookmarks("TmpSplit1").Delete()

$p2 = $d.Content
$p2.Find.Execute("PROBLEMI:Reload, ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$split2 = $p2.End
$d.Bookmarks.Add("TmpSplit2", $d.Range($split2, $split2))
$d.Bookmarks("TmpSplit2").Delete()
